$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------------
# 1) Re-add the (reverted) Russian localisation duplicate-source-text lines.
#    These two shared strings previously held only the English source line;
#    restore the "<english>\n---\n<english>" pattern used elsewhere in the
#    sheet for rows that still need RU translation.
# -------------------------------------------------------------------------
$warning = "Warning commander, sensors are reading multiple contacts approaching from multiple vectors. Looks like the attacking force has brought reinforcements."
$ws.Range("B41").Value = $warning + "`n---`n" + $warning + " "

$noNoNo = "No, No, No. These things are old and fragile, they can't be shaken like this."
$ws.Range("B42").Value = $noNoNo + "`n---`n" + $noNoNo + " "

# Embedded newlines above make Excel auto-expand the row height; re-fit it
# back down so the row keeps its original (default) height.
$ws.Rows("41:42").AutoFit()

# -------------------------------------------------------------------------
# 2) Fix the double space typo in the English source text.
# -------------------------------------------------------------------------
$ws.Range("B111").Value = "Destroy the Second Pirate"

# -------------------------------------------------------------------------
# 3) Move the "needs attention" red highlight: it used to mark the C column
#    of rows 78-81; now it marks the whole of rows 41-42 instead.
# -------------------------------------------------------------------------
$ws.Range("C78:C81").Interior.Pattern = -4142
$ws.Range("A41:C42").Interior.Color = 255

# -------------------------------------------------------------------------
# 4) Reset the saved scroll position back to the top of the sheet.
# -------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# -------------------------------------------------------------------------
# 5) Column B used to share a single run with the rest of the sheet; touch
#    its width explicitly so it becomes its own <col> entry again.
# -------------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = $ws.Columns("B").ColumnWidth
